$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E are treated as text so numeric-looking values
# (e.g. "54.298.79", "1.00", "239.00") are preserved exactly as strings
# instead of being coerced into numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '54.298.79', '  +0.48%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.277.11', '  +0.72%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.08%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '497.91', '  +1.32%  ')
    ,@(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '127.95', '  +0.79%  ')
    ,@(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.999', '  -0.15%  ')
    ,@(8, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.528', '  +0.29%  ')
    ,@(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.0955', '  +2.54%  ')
    ,@(10, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.152', '  +1.43%  ')
    ,@(11, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.333', '  +3.01%  ')
    ,@(12, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '4.70', '  +1.49%  ')
    ,@(13, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.684.13', '  +0.89%  ')
    ,@(14, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '22.55', '  +5.13%  ')
    ,@(15, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '54.249.75', '  +0.40%  ')
    ,@(16, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000129', '  +0.26%  ')
    ,@(17, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.294.40', '  +0.69%  ')
    ,@(18, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '10.25', '  +4.90%  ')
    ,@(19, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.12', '  +2.15%  ')
    ,@(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '304.75', '  +2.41%  ')
    ,@(21, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.42', '  +2.81%  ')
    ,@(22, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.999', '  -0.07%  ')
    ,@(23, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '61.94', '  -2.93%  ')
    ,@(24, 'Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '0.999', '  -0.20%  ')
    ,@(25, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.151', '  +2.76%  ')
    ,@(26, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '7.33', '  +2.97%  ')
    ,@(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '175.18', '  +7.64%  ')
    ,@(28, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.61', '  +1.20%  ')
    ,@(29, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '5.94', '  +2.55%  ')
    ,@(30, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0684', '  +0.86%  ')
    ,@(31, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '1.08', '  +1.46%  ')
    ,@(32, 'USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '0.999', '  -0.01%  ')
    ,@(33, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.75', '  +1.93%  ')
    ,@(34, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.998', '  +0.02%  ')
    ,@(35, 'SuiNetwork', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui', '0.925', '  +10.34%  ')
    ,@(36, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.20', '  +0.80%  ')
    ,@(37, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '3.74', '  +3.26%  ')
    ,@(38, 'PolygonEcosystemToken', 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol', '0.373', '  -0.30%  ')
    ,@(39, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.41', '  +1.41%  ')
    ,@(40, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.38', '  +1.56%  ')
    ,@(41, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '125.91', '  -0.58%  ')
    ,@(42, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '4.76', '  -0.60%  ')
    ,@(43, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0490', '  +2.57%  ')
    ,@(44, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.0896', '  +0.52%  ')
    ,@(45, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.545', '  +0.47%  ')
    ,@(46, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '239.00', '  -1.29%  ')
    ,@(47, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.372', '  -0.42%  ')
    ,@(48, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0205', '  +1.38%  ')
    ,@(49, 'WhiteBITCoin', 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt', '10.77', '  +0.95%  ')
    ,@(50, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '16.27', '  +0.49%  ')
    ,@(51, 'ZEEBU', 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu', '4.64', '  +0.46%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
